$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First, copy the existing row 2 (original data) down into new row 3,
# so row 3 preserves the original (pre-edit) values that were in row 2.
$ws.Rows.Item(2).Copy()
$ws.Rows.Item(3).PasteSpecial(-4104)  # xlPasteAll
$excel.CutCopyMode = 0

# --- Now update row 2 with the new record's data ---

# New record id
$ws.Range("A2").Value = 81397446

# Updated location name (a shortened version, without "Näsum, ")
$ws.Range("P2").Value = "Hejabacken, Sk"

# Updated coordinates
$ws.Range("Q2").Value = 465267.8687842482
$ws.Range("R2").Value = 6222143.925440002

# Updated start/end dates -- keep as plain text (not an Excel date serial)
$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value = "2008-05-02"
$ws.Range("Y2").Style = "Normal"

$ws.Range("AA2").NumberFormat = "@"
$ws.Range("AA2").Value = "2008-05-02"
$ws.Range("AA2").Style = "Normal"

# New habitat description value
$ws.Range("AI2").Value = "Bokskog"

# Updated observers
$ws.Range("AX2").Value = "Krister Håkansson"
